$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 2
$cols = @(2,3,4,5,6,9,10,11,12,13,14)
$data = @(
    @(1.02,1.051273764027869,1.05118370490395,0.992614727750844,1.058976453859706,1.048681408022068,1.056302504192728,1.053935912474021,0.9955398523335997,1.061707217612974,1.057802574908272),
    @(1.02,1.052533998183304,1.052156091355224,0.9936372048519299,1.060441499866155,1.049143712843775,1.05721124091216,1.054720483195767,0.9963617723202687,1.062984754820608,1.058712602138079),
    @(1.02,1.053347991943685,1.052783875031494,0.9942998659930998,1.06138847480587,1.049440737737445,1.057797289167355,1.055226123099598,0.9968940712668347,1.063809885846709,1.059299482649018),
    @(1.02,1.053689849633786,1.053047459533702,0.994578699834602,1.061786347606994,1.049565102189247,1.058043196879179,1.055438211191599,0.9971179600053012,1.064156412002906,1.059545739577999),
    @(1.02,1.053747228933189,1.05309169693204,0.994625531979634,1.061853138561521,1.049585953985012,1.058084458583507,1.055473793484025,0.9971555583673455,1.06421457437224,1.059587059878679),
    @(1.02,1.05335256121375,1.052787398375972,0.994303590798249,1.061393792119713,1.049442401482363,1.057800576828801,1.055228958924882,0.9968970624462089,1.063814517550058,1.05930277497932),
    @(1.02,1.051699971627777,1.051512622386632,0.9929600610674297,1.059471785187632,1.048838085731696,1.056610025029351,1.054201483931317,0.9958175282591056,1.06213928535659,1.058110532459756),
    @(1.02,1.048776488147909,1.049255324891316,0.9906006454969559,1.056076984236263,1.047756895953692,1.054496912710585,1.052375259874171,0.9939188001724441,1.059175429826549,1.055994419279121),
    @(1.02,1.046819528998376,1.047742880924308,0.989033133672735,1.053808052307852,1.047025006986666,1.053077724416754,1.051147041628869,0.9926553831429383,1.057191208536808,1.054573215575128),
    @(1.02,1.045970186420443,1.047086133755474,0.988355674866747,1.052824138169953,1.046705429590609,1.052460672379776,1.050612619917919,0.9921088820399291,1.056329972254818,1.053955287253481),
    @(1.02,1.045654401458129,1.046841906891274,0.9881042295826724,1.052458443393187,1.046586321358274,1.052231086649936,1.05041371819107,0.9919059725120875,1.056009755515142,1.053725375485571),
    @(1.02,1.045722152134489,1.04689430719488,0.9881581567098651,1.05253689647308,1.046611888755496,1.05228035106726,1.050456401164461,0.9919494934313052,1.056078457461871,1.053774709864019),
    @(1.02,1.045944089717826,1.047065951656672,0.9883348863814464,1.052793914362412,1.046695592302447,1.052441702633821,1.050596186692866,0.9920921077337197,1.056303509507395,1.053936290568311),
    @(1.02,1.046080792684183,1.047171670025028,0.9884438009545853,1.052952241502815,1.046747111345554,1.052541065493677,1.050682260885446,0.9921799884222134,1.056442129534781,1.054035794534826),
    @(1.02,1.046875855137422,1.047786427780059,0.9890781214508737,1.053873320253861,1.047046159921259,1.053118622361675,1.051182454436139,0.9926916645766087,1.057248322099254,1.054614171599823),
    @(1.02,1.047374046214513,1.048171551288757,0.9894763578477731,1.054450694930388,1.047233029943086,1.053480226965446,1.051495515160988,0.9930127773692701,1.057753470579613,1.054976289723608),
    @(1.02,1.047664443317239,1.04839600918719,0.9897087662937551,1.054787328235653,1.047341771144515,1.053690900604851,1.051677868055397,0.9932001317071766,1.058047917266423,1.055187262543745),
    @(1.02,1.047763429289375,1.048472513434792,0.9897880325774039,1.054902088034713,1.047378805584835,1.053762693582109,1.051740003368596,0.993264023964098,1.058148282608282,1.055259157475266),
    @(1.02,1.047320614653563,1.048130249635369,0.9894336180360677,1.054388762582308,1.04721300715251,1.053441455506296,1.051461952608812,0.9929783193494215,1.057699293443725,1.054937463204539),
    @(1.02,1.04587874299636,1.047015414449552,0.9882828385668249,1.052718235275976,1.046670954838924,1.052394199286012,1.050555034237922,0.9920501090198102,1.056237246043319,1.0538887197603),
    @(1.02,1.044970434815412,1.046312841379829,0.9875604150241495,1.051666601605748,1.046327812329876,1.051733517294475,1.049982538876018,0.9914670000341481,1.055316171724123,1.053227099524554),
    @(1.02,1.045452113667264,1.046685444668867,0.9879432794636464,1.05222421878938,1.046509940660904,1.052083970216066,1.050286246879399,0.9917760702887611,1.055804625995961,1.053578050129491),
    @(1.02,1.04734475866062,1.048148912624901,0.9894529299347244,1.054416747575323,1.047222055383373,1.053458975414071,1.051477118854081,0.9929938892766441,1.057723774370356,1.054955007992592),
    @(1.02,1.049533661458441,1.049840211707189,0.9912096547607051,1.056955604830508,1.048038355470513,1.055045027854726,1.052849260627809,0.9944092447426416,1.059943099349109,1.056543312809636)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($row, $cols[$j]).Value = $rowData[$j]
    }
}
